$d = $word.ActiveDocument

function Add-Para([string]$text) {
    $r = $d.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $r.Collapse(0)
    if ($text -ne "") {
        $r.Text = $text
    }
    return $d.Paragraphs.Item($d.Paragraphs.Count)
}

# 1) Empty paragraph with a bottom border (sz=24, default font) -------------
$p1 = Add-Para ""
$brd1 = $p1.Borders.Item(-3)
$brd1.LineStyle = 1
$brd1.LineWidth = 2
$brd1.Color = -16777216
$p1.Borders.DistanceFromBottom = 0

# 2) Empty paragraph, no border, sz=24, default font -------------------------
$p2 = Add-Para ""
$brd2 = $p2.Borders.Item(-3)
$brd2.LineStyle = 0

# 3) "Personal OKR (update)" - bold, Arial, sz=22 -----------------------------
$p3 = Add-Para "Personal OKR (update)"
$p3.Range.Font.Name = "Arial"
$p3.Range.Font.NameBi = "Arial"
$p3.Range.Font.Size = 11
$p3.Range.Font.SizeBi = 11
$p3.Range.Font.Bold = $true
$p3.Range.Font.BoldBi = $true

# Move the _GoBack bookmark to the end of this paragraph's text (collapsed).
$bookmarkPos = $p3.Range.End - 1
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 4) Empty paragraph, Arial, sz=22, not bold ---------------------------------
$p4 = Add-Para ""
$p4.Range.Font.Bold = $false
$p4.Range.Font.BoldBi = $false

# 5-10) Bullet-style lines, Arial, sz=22, not bold ---------------------------
Add-Para "*Create a financially healty and stable life" | Out-Null
Add-Para "*more workouts and get into shape" | Out-Null
Add-Para "*Create a portfolio I can be proud of" | Out-Null
Add-Para "*Increase art skills and combine with coding (2d + 3D )" | Out-Null
Add-Para "*Work out concepts for apps + games and wrhite USP documents for each one." | Out-Null
Add-Para "*research potential employees" | Out-Null

# 11) Empty paragraph, bold, Arial, sz=22 ------------------------------------
$p11 = Add-Para ""
$p11.Range.Font.Bold = $true
$p11.Range.Font.BoldBi = $true

# 12) "Pro OKR" bold, Arial, sz=22 --------------------------------------------
Add-Para "Pro OKR" | Out-Null

# 13-14) Bullet-style lines, Arial, sz=22, not bold ---------------------------
$p13 = Add-Para "*Update linkedin acount with photo"
$p13.Range.Font.Bold = $false
$p13.Range.Font.BoldBi = $false
Add-Para "*post more often on Instagram and Linkedin and stay connected" | Out-Null

# 15) Empty paragraph, Arial, sz=22, not bold ---------------------------------
Add-Para "" | Out-Null

# 16) Empty paragraph, default font, sz=24 ------------------------------------
$p16 = Add-Para ""
$p16.Range.Font.Name = "Calibri"
$p16.Range.Font.NameBi = "Calibri"
$p16.Range.Font.Size = 12
$p16.Range.Font.SizeBi = 12

# Update the doc-defaults run fonts (rPrDefault) from theme fonts to explicit
# Times New Roman / SimSun, matching styles.xml's docDefaults change.
$normalStyle = $d.Styles.Item(-1)
$normalStyle.Font.Name = "Times New Roman"
$normalStyle.Font.NameFarEast = "SimSun"
$normalStyle.Font.NameBi = "Times New Roman"
